# Append/update the "取得日時" (retrieved-at) timestamp for the newly
# scraped rows (2-6) on the active sheet ("ランサーズ") to reflect the
# latest scrape run time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-12-07 06:33:09"

for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
